# New submission synced into the "JSS 3E" response sheet:
#  - C2 (Admission No for the existing row) was entered as text "15" and is
#    corrected/recognised as the number 15.
#  - A brand-new form response lands in row 3.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3E")

# Admission No for row 2 becomes a real number.
$ws.Cells.Item(2, 3).Value = 15

# New submission appended as row 3.
$ws.Cells.Item(3, 1).Value = "2026-02-07 19:50:53"
$ws.Cells.Item(3, 2).Value = "Moduye Khadija "
# Admission No "45" is kept as text for this response (leading apostrophe
# forces Excel to store the digits as a literal string, not a number).
$ws.Cells.Item(3, 3).Value = "'45"
$ws.Cells.Item(3, 4).Value = 9
